$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.278738260269165
$ws.Range("B1").Value = 1.021008729934692
$ws.Range("C1").Value = 4.255475521087646
$ws.Range("D1").Value = 2.554011106491089
$ws.Range("E1").Value = 0.7529987096786499
